$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.858.24'
$ws.Range("E2").Value = '  +1.71%  '

# Row 3
$ws.Range("D3").Value = '1.872.50'
$ws.Range("E3").Value = '  +1.91%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.032'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.25'
$ws.Range("E5").Value = '  +2.03%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.027'
$ws.Range("E6").Value = '  +0.35%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4435'
$ws.Range("E7").Value = '  +1.72%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3826'
$ws.Range("E8").Value = '  +2.94%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07468'
$ws.Range("E9").Value = '  +1.69%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8903'
$ws.Range("E10").Value = '  +2.28%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.79'
$ws.Range("E11").Value = '  +2.26%  '

# Row 12
$ws.Range("D12").Value = '1.871.82'
$ws.Range("E12").Value = '  -3.66%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.588'
$ws.Range("E13").Value = '  +2.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.787'
$ws.Range("E14").Value = '  +1.57%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07193'
$ws.Range("E15").Value = '  +1.30%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '85.84'
$ws.Range("E16").Value = '  +4.45%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.033'
$ws.Range("E17").Value = '  +0.45%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009145'
$ws.Range("E18").Value = '  +1.89%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.028'
$ws.Range("E19").Value = '  +0.50%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.62'
$ws.Range("E20").Value = '  +1.62%  '

# Row 21
$ws.Range("D21").Value = '27.882.20'
$ws.Range("E21").Value = '  +1.68%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.326'
$ws.Range("E22").Value = '  +1.56%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.34'
$ws.Range("E23").Value = '  +1.76%  '

# Row 24
$ws.Range("D24").Value = '2.095.99'
$ws.Range("E24").Value = '  -1.92%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.025'
$ws.Range("E25").Value = '  +6.21%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.51'
$ws.Range("E26").Value = '  +1.14%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.91'
$ws.Range("E27").Value = '  +2.05%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.439'
$ws.Range("E28").Value = '  +3.86%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.994'
$ws.Range("E29").Value = '  +3.92%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.53'
$ws.Range("E30").Value = '  +2.59%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09075'
$ws.Range("E31").Value = '  +0.44%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.236'
$ws.Range("E32").Value = '  +3.23%  '

# Row 33
$ws.Range("E33").Value = '  +3.61%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.619'
$ws.Range("E34").Value = '  +3.60%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.016'
$ws.Range("E35").Value = '  +5.37%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.030'
$ws.Range("E36").Value = '  +0.41%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.149'
$ws.Range("E37").Value = '  +0.31%  '

# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05371'
$ws.Range("E38").Value = '  +2.44%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01990'
$ws.Range("E39").Value = '  +1.85%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.886'
$ws.Range("E40").Value = '  +3.46%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5235'
$ws.Range("E41").Value = '  +1.56%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1699'
$ws.Range("E42").Value = '  +2.40%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.939'
$ws.Range("E43").Value = '  +6.04%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.844'
$ws.Range("E44").Value = '  +4.52%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '112.46'
$ws.Range("E45").Value = '  +3.82%  '

# Row 46
$ws.Range("E46").Value = '  +2.43%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.06622'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.726'
$ws.Range("E48").Value = '  +3.12%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.029'
$ws.Range("E49").Value = '  +0.33%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4759'
$ws.Range("E50").Value = '  +3.08%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.930'
$ws.Range("E51").Value = '  +2.89%  '
